$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11241334
$ws.Range("I32").Value = 11910194
$ws.Range("K32").Value = 11910194
$ws.Range("M32").Value = -11909907
$ws.Range("H45").Value = 2414.8
$ws.Range("J45").Value = 2904.6667
$ws.Range("L45").Value = 2904.6667
$ws.Range("N45").Value = -3658.6667
$ws.Range("H74").Value = 2545.442
$ws.Range("I74").Value = 2723.182
$ws.Range("J74").Value = 1958.9
$ws.Range("K74").Value = 2723.182
$ws.Range("L74").Value = 1958.9
$ws.Range("M74").Value = -1849.182
$ws.Range("N74").Value = -3706.9
$ws.Range("H76").Value = 68166.664
$ws.Range("J76").Value = 68166.664
$ws.Range("L76").Value = 68166.664
$ws.Range("N76").Value = -68842.664
$ws.Range("H77").Value = 2545.442
$ws.Range("I77").Value = 2723.182
$ws.Range("J77").Value = 1958.9
$ws.Range("K77").Value = 13615.91
$ws.Range("L77").Value = 9794.5
$ws.Range("M77").Value = -9247.91
$ws.Range("N77").Value = -18530.5
$ws.Range("H79").Value = 68166.664
$ws.Range("J79").Value = 68166.664
$ws.Range("L79").Value = 68166.664
$ws.Range("N79").Value = -70506.664
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000
$ws.Range("N95").Value = -105492
$ws.Range("H102").Value = 1542.421
$ws.Range("I102").Value = 1544.1666
$ws.Range("K102").Value = 1544.1666
$ws.Range("M102").Value = 77.83339999999998
$ws.Range("H132").Value = 2464.4285
$ws.Range("I132").Value = 2208.5557
$ws.Range("K132").Value = 6625.6671
$ws.Range("M132").Value = -4095.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1307.1428
$ws.Range("I20").Value = 1348.0769
$ws.Range("J20").Value = 1240.625
$ws.Range("K20").Value = 1348.0769
$ws.Range("L20").Value = 1240.625
$ws.Range("M20").Value = -1101.0769
$ws.Range("N20").Value = -1734.625
$ws.Range("H100").Value = 17000
$ws.Range("J100").Value = 17000
$ws.Range("L100").Value = 17000
$ws.Range("N100").Value = -19164
$ws.Range("H105").Value = 2977.0688
$ws.Range("I105").Value = 1255.3529
$ws.Range("K105").Value = 1255.3529
$ws.Range("M105").Value = 491.6470999999999
$ws.Range("H134").Value = 1621.5964
$ws.Range("I134").Value = 1175.1177
$ws.Range("K134").Value = 3525.3531
$ws.Range("M134").Value = -990.3531000000003
$ws.Range("H141").Value = 44995
$ws.Range("J141").Value = 44995
$ws.Range("L141").Value = 44995
$ws.Range("N141").Value = -55355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2139.4634
$ws.Range("I31").Value = 1883.1818
$ws.Range("J31").Value = 3196.625
$ws.Range("K31").Value = 1883.1818
$ws.Range("L31").Value = 3196.625
$ws.Range("M31").Value = -1588.1818
$ws.Range("N31").Value = -3786.625
$ws.Range("H34").Value = 2139.4634
$ws.Range("I34").Value = 1883.1818
$ws.Range("J34").Value = 3196.625
$ws.Range("K34").Value = 1883.1818
$ws.Range("L34").Value = 3196.625
$ws.Range("M34").Value = -1681.1818
$ws.Range("N34").Value = -3600.625
$ws.Range("H82").Value = 24000
$ws.Range("J82").Value = 24000
$ws.Range("L82").Value = 24000
$ws.Range("N82").Value = -24722
$ws.Range("H85").Value = 24000
$ws.Range("J85").Value = 24000
$ws.Range("L85").Value = 24000
$ws.Range("N85").Value = -26496
$ws.Range("H86").Value = 50338.89
$ws.Range("I86").Value = 66995.2
$ws.Range("J86").Value = 29518.5
$ws.Range("K86").Value = 66995.2
$ws.Range("L86").Value = 29518.5
$ws.Range("M86").Value = -65872.2
$ws.Range("N86").Value = -31764.5
$ws.Range("H89").Value = 50338.89
$ws.Range("I89").Value = 66995.2
$ws.Range("J89").Value = 29518.5
$ws.Range("K89").Value = 334976
$ws.Range("L89").Value = 147592.5
$ws.Range("M89").Value = -329360
$ws.Range("N89").Value = -158824.5
$ws.Range("H105").Value = 4927.615
$ws.Range("I105").Value = 4414.4546
$ws.Range("K105").Value = 4414.4546
$ws.Range("M105").Value = -2667.4546
$ws.Range("H122").Value = 734685.1
$ws.Range("I122").Value = 1703598.6
$ws.Range("K122").Value = 5110795.800000001
$ws.Range("M122").Value = -5108345.800000001
$ws.Range("H132").Value = 1881.6774
$ws.Range("I132").Value = 1776.9656
$ws.Range("K132").Value = 5330.8968
$ws.Range("M132").Value = -2800.8968
$ws.Range("H134").Value = 3189.258
$ws.Range("I134").Value = 2708.5264
$ws.Range("J134").Value = 3950.4167
$ws.Range("K134").Value = 8125.5792
$ws.Range("L134").Value = 11851.2501
$ws.Range("M134").Value = -5590.5792
$ws.Range("N134").Value = -16921.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 717.6667
$ws.Range("I68").Value = 632.6
$ws.Range("J68").Value = 824
$ws.Range("K68").Value = 1897.8
$ws.Range("L68").Value = 2472
$ws.Range("M68").Value = -1086.8
$ws.Range("N68").Value = -4094
$ws.Range("H71").Value = 717.6667
$ws.Range("I71").Value = 632.6
$ws.Range("J71").Value = 824
$ws.Range("K71").Value = 5693.400000000001
$ws.Range("L71").Value = 7416
$ws.Range("M71").Value = -1637.400000000001
$ws.Range("N71").Value = -15528
$ws.Range("H107").Value = 263
$ws.Range("I107").Value = 295
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 885
$ws.Range("L107").Value = 597
$ws.Range("M107").Value = 1035
$ws.Range("N107").Value = -4437
$ws.Range("H136").Value = 1692.1666
$ws.Range("I136").Value = 1030.6
$ws.Range("K136").Value = 3091.8
$ws.Range("M136").Value = 2008.2
$ws.Range("H137").Value = 3241.9473
$ws.Range("I137").Value = 2812.5
$ws.Range("J137").Value = 3554.2727
$ws.Range("K137").Value = 8437.5
$ws.Range("L137").Value = 10662.8181
$ws.Range("M137").Value = -3337.5
$ws.Range("N137").Value = -20862.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 833
$ws.Range("J40").Value = 999.5
$ws.Range("L40").Value = 999.5
$ws.Range("N40").Value = -1301.5
$ws.Range("H107").Value = 86.59999999999999
$ws.Range("J107").Value = 90
$ws.Range("L107").Value = 90
$ws.Range("N107").Value = -3930
$ws.Range("H126").Value = 3227.4443
$ws.Range("I126").Value = 3227.4443
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9682.332900000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7212.332900000001
$ws.Range("N126").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1528.2106
$ws.Range("I22").Value = 1299.8
$ws.Range("J22").Value = 1782
$ws.Range("K22").Value = 1299.8
$ws.Range("L22").Value = 1782
$ws.Range("M22").Value = -1004.8
$ws.Range("N22").Value = -2372
$ws.Range("H27").Value = 1528.2106
$ws.Range("I27").Value = 1299.8
$ws.Range("J27").Value = 1782
$ws.Range("K27").Value = 1299.8
$ws.Range("L27").Value = 1782
$ws.Range("M27").Value = -1192.8
$ws.Range("N27").Value = -1996
$ws.Range("H48").Value = 7646
$ws.Range("J48").Value = 7646
$ws.Range("L48").Value = 7646
$ws.Range("N48").Value = -8968
$ws.Range("H93").Value = 3191.3635
$ws.Range("I93").Value = 3498.6
$ws.Range("J93").Value = 119
$ws.Range("K93").Value = 3498.6
$ws.Range("L93").Value = 119
$ws.Range("M93").Value = -2250.6
$ws.Range("N93").Value = -2615
$ws.Range("H132").Value = 3217.327
$ws.Range("I132").Value = 2085.5
$ws.Range("K132").Value = 6256.5
$ws.Range("M132").Value = -3726.5
$ws.Range("H136").Value = 1885.4222
$ws.Range("I136").Value = 1758.1904
$ws.Range("K136").Value = 5274.5712
$ws.Range("M136").Value = -2724.5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 50000
$ws.Range("I26").Value = 50000
$ws.Range("K26").Value = 50000
$ws.Range("M26").Value = -49707
$ws.Range("H29").Value = 1200
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
$ws.Range("H122").Value = 2080.7896
$ws.Range("I122").Value = 2205.2415
$ws.Range("J122").Value = 1679.7778
$ws.Range("K122").Value = 6615.7245
$ws.Range("L122").Value = 5039.3334
$ws.Range("M122").Value = -4165.7245
$ws.Range("N122").Value = -9939.3334
$ws.Range("H126").Value = 1625.6923
$ws.Range("I126").Value = 1460.4445
$ws.Range("K126").Value = 4381.333500000001
$ws.Range("M126").Value = -1911.333500000001
$ws.Range("H136").Value = 1380.9535
$ws.Range("I136").Value = 860.97144
$ws.Range("K136").Value = 2582.91432
$ws.Range("M136").Value = -32.91431999999986
